$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "SUM" total row (added first so the new shared strings land in the
# same order as the target file: SUM, (+1), Giulio A. Abbo)
$ws.Range("A12").Value = "SUM"
$ws.Range("A12").HorizontalAlignment = -4152   # xlRight
$ws.Range("B12").Formula = "=SUM(B3:B11)"

# Mark the two tasks that got an extra hour of effort with "(+1)"
$ws.Range("C3").Value = "(+1)"
$ws.Range("C4").Value = "(+1)"

# Author name tweak
$ws.Range("A1").Value = "Giulio A. Abbo"

# Updated effort hours per task
$ws.Range("B4").Value = 3      # Product perspective: 3.5 -> 3
$ws.Range("B5").Value = 2.5    # Product functions: 2 -> 2.5
$ws.Range("B6").Value = 1      # Domain assumptions: 0.75 -> 1
$ws.Range("B7").Value = 0      # External interface requirements: (blank) -> 0
$ws.Range("B9").Value = 0      # Non-functional Requirements: (blank) -> 0
$ws.Range("B10").Value = 0     # Formal analysis using Alloy: (blank) -> 0

# Restore the active selection to B4 and set print/page setup
$ws.Range("B4").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
